$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "kelas" (column D) and "jurusan" (column E) for each student row.
# Rows 2-15 -> kelas 1027404, rows 16-26 -> kelas 1028404.
# All rows 2-26 -> jurusan "IPA" (replacing "TJKT").
for ($r = 2; $r -le 26; $r++) {
    if ($r -le 15) {
        $ws.Cells.Item($r, 4).Value = 1027404
    } else {
        $ws.Cells.Item($r, 4).Value = 1028404
    }
    $ws.Cells.Item($r, 5).Value = "IPA"
}

# Update the active selection shown in the worksheet view.
$ws.Range("D16:D26").Select() | Out-Null
